$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores": update the pairwise p-value matrix (B2:F6) ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8804346223621609
$wsP.Range("D2").Value = 0.8620072575740449
$wsP.Range("E2").Value = 0.7570713489560053
$wsP.Range("F2").Value = 0.6119780666839594

$wsP.Range("B3").Value = 0.8804346223621609
$wsP.Range("D3").Value = 0.9932350104831833
$wsP.Range("E3").Value = 0.8910507358777964
$wsP.Range("F3").Value = 0.548277384886831

$wsP.Range("B4").Value = 0.8620072575740449
$wsP.Range("C4").Value = 0.9932350104831833
$wsP.Range("E4").Value = 0.8140666050637364
$wsP.Range("F4").Value = 0.5760571970857356

$wsP.Range("B5").Value = 0.7570713489560053
$wsP.Range("C5").Value = 0.8910507358777964
$wsP.Range("D5").Value = 0.8140666050637364
$wsP.Range("F5").Value = 0.6972358966290013

$wsP.Range("B6").Value = 0.6119780666839594
$wsP.Range("C6").Value = 0.548277384886831
$wsP.Range("D6").Value = 0.5760571970857356
$wsP.Range("E6").Value = 0.6972358966290013

# --- Sheet "Estadisticos_DM": update the pairwise DM statistic matrix (B2:F6) ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -0.1531898747556643
$wsE.Range("D2").Value = -0.1770470798264044
$wsE.Range("E2").Value = -0.3154524618924585
$wsE.Range("F2").Value = -0.5188368313714822

$wsE.Range("B3").Value = 0.1531898747556643
$wsE.Range("D3").Value = -0.008631406053872769
$wsE.Range("E3").Value = -0.1394887462470434
$wsE.Range("F3").Value = -0.6152242873740289

$wsE.Range("B4").Value = 0.1770470798264044
$wsE.Range("C4").Value = 0.008631406053872769
$wsE.Range("E4").Value = -0.2396606877006016
$wsE.Range("F4").Value = -0.5725078924192584

$wsE.Range("B5").Value = 0.3154524618924585
$wsE.Range("C5").Value = 0.1394887462470434
$wsE.Range("D5").Value = 0.2396606877006016
$wsE.Range("F5").Value = -0.397163346971767

$wsE.Range("B6").Value = 0.5188368313714822
$wsE.Range("C6").Value = 0.6152242873740289
$wsE.Range("D6").Value = 0.5725078924192584
$wsE.Range("E6").Value = 0.397163346971767
